# Updates Sheets per scheduled runner diff: refresh market-price derived columns
# (currentAveragePrice, currentAveragePriceNQ/HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 98
$ws.Range("H98").Value = 796.0769
$ws.Range("I98").Value = 762.5
$ws.Range("K98").Value = 762.5
$ws.Range("M98").Value = 735.5
# Row 112
$ws.Range("H112").Value = 27101948
$ws.Range("I112").Value = 250000510
$ws.Range("J112").Value = 3004806.8
$ws.Range("K112").Value = 750001530
$ws.Range("L112").Value = 9014420.399999999
$ws.Range("M112").Value = -750000422
$ws.Range("N112").Value = -9016636.399999999
# Row 122
$ws.Range("H122").Value = 796.0769
$ws.Range("I122").Value = 762.5
$ws.Range("K122").Value = 2287.5
$ws.Range("M122").Value = 162.5
# Row 129
$ws.Range("H129").Value = 871.40625
$ws.Range("J129").Value = 1007.73914
$ws.Range("L129").Value = 3023.21742
$ws.Range("N129").Value = -13023.21742

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 5
$ws.Range("H5").Value = 7519063.5
$ws.Range("J5").Value = 1000
$ws.Range("L5").Value = 1000
$ws.Range("N5").Value = -1224
# Row 61
$ws.Range("H61").Value = 55667940
$ws.Range("I61").Value = 71500850
$ws.Range("J61").Value = 252753.5
$ws.Range("K61").Value = 71500850
$ws.Range("L61").Value = 252753.5
$ws.Range("M61").Value = -71500638
$ws.Range("N61").Value = -253177.5
# Row 74
$ws.Range("H74").Value = 9002567
$ws.Range("I74").Value = 13212514
$ws.Range("J74").Value = 114900
$ws.Range("K74").Value = 13212514
$ws.Range("L74").Value = 114900
$ws.Range("M74").Value = -13211640
$ws.Range("N74").Value = -116648
# Row 77
$ws.Range("H77").Value = 9002567
$ws.Range("I77").Value = 13212514
$ws.Range("J77").Value = 114900
$ws.Range("K77").Value = 66062570
$ws.Range("L77").Value = 574500
$ws.Range("M77").Value = -66058202
$ws.Range("N77").Value = -583236
# Row 122
$ws.Range("H122").Value = 3833341.2
$ws.Range("I122").Value = 1935.8
$ws.Range("J122").Value = 27779624
$ws.Range("K122").Value = 5807.4
$ws.Range("L122").Value = 83338872
$ws.Range("M122").Value = -3357.4
$ws.Range("N122").Value = -83343772
# Row 136
$ws.Range("H136").Value = 55667940
$ws.Range("I136").Value = 71500850
$ws.Range("J136").Value = 252753.5
$ws.Range("K136").Value = 214502550
$ws.Range("L136").Value = 758260.5
$ws.Range("M136").Value = -214500000
$ws.Range("N136").Value = -763360.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("H4").Value = 7519063.5
$ws.Range("J4").Value = 1000
$ws.Range("L4").Value = 1000
$ws.Range("N4").Value = -1230
# Row 13
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
# Row 20
$ws.Range("H20").Value = 1586.1666
$ws.Range("I20").Value = 1104
$ws.Range("J20").Value = 1827.25
$ws.Range("K20").Value = 1104
$ws.Range("L20").Value = 1827.25
$ws.Range("M20").Value = -857
$ws.Range("N20").Value = -2321.25
# Row 94
$ws.Range("H94").Value = 1100
$ws.Range("I94").Value = 1300
$ws.Range("J94").Value = 1000
$ws.Range("K94").Value = 1300
$ws.Range("L94").Value = 1000
$ws.Range("M94").Value = -849
$ws.Range("N94").Value = -1902
# Row 99
$ws.Range("H99").Value = 1104.6364
$ws.Range("I99").Value = 1169.5238
$ws.Range("J99").Value = 991.0833
$ws.Range("K99").Value = 1169.5238
$ws.Range("L99").Value = 991.0833
$ws.Range("M99").Value = 328.4762000000001
$ws.Range("N99").Value = -3987.0833
# Row 134
$ws.Range("H134").Value = 2839.1333
$ws.Range("I134").Value = 2271.2222
$ws.Range("K134").Value = 6813.6666
$ws.Range("M134").Value = -4278.6666

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1518.75
$ws.Range("I16").Value = 1150
$ws.Range("J16").Value = 1740
$ws.Range("K16").Value = 1150
$ws.Range("L16").Value = 1740
$ws.Range("M16").Value = -863
$ws.Range("N16").Value = -2314
# Row 19
$ws.Range("H19").Value = 679.8
$ws.Range("I19").Value = 350
$ws.Range("K19").Value = 350
$ws.Range("M19").Value = -180
# Row 24
$ws.Range("H24").Value = 679.8
$ws.Range("I24").Value = 350
$ws.Range("K24").Value = 350
$ws.Range("M24").Value = -180
# Row 94
$ws.Range("H94").Value = 4605.4614
$ws.Range("I94").Value = 12102.75
$ws.Range("K94").Value = 12102.75
$ws.Range("M94").Value = -11651.75
# Row 113
$ws.Range("H113").Value = 1518.75
$ws.Range("I113").Value = 1150
$ws.Range("J113").Value = 1740
$ws.Range("K113").Value = 1150
$ws.Range("L113").Value = 1740
$ws.Range("M113").Value = 1020
$ws.Range("N113").Value = -6080
# Row 132
$ws.Range("H132").Value = 28002.553
$ws.Range("I132").Value = 1488.5428
$ws.Range("J132").Value = 337332.66
$ws.Range("K132").Value = 4465.6284
$ws.Range("L132").Value = 1011997.98
$ws.Range("M132").Value = -1935.6284
$ws.Range("N132").Value = -1017057.98

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 15
$ws.Range("H15").Value = 516
$ws.Range("I15").Value = 100
$ws.Range("J15").Value = 793.3333
$ws.Range("K15").Value = 300
$ws.Range("L15").Value = 2379.9999
$ws.Range("M15").Value = -160
$ws.Range("N15").Value = -2659.9999

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Range("H122").Value = 3296.9412
$ws.Range("I122").Value = 2860.5715
$ws.Range("J122").Value = 5333.3335
$ws.Range("K122").Value = 8581.7145
$ws.Range("L122").Value = 16000.0005
$ws.Range("M122").Value = -6131.7145
$ws.Range("N122").Value = -20900.0005
# Row 126
$ws.Range("H126").Value = 1940
$ws.Range("I126").Value = 1800
$ws.Range("J126").Value = 2266.6667
$ws.Range("K126").Value = 5400
$ws.Range("L126").Value = 6800.000100000001
$ws.Range("M126").Value = -2930
$ws.Range("N126").Value = -11740.0001
# Row 132
$ws.Range("H132").Value = 55878.51
$ws.Range("I132").Value = 39867.04
$ws.Range("J132").Value = 93723.82000000001
$ws.Range("K132").Value = 119601.12
$ws.Range("L132").Value = 281171.46
$ws.Range("M132").Value = -117071.12
$ws.Range("N132").Value = -286231.46

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2188.7083
$ws.Range("I7").Value = 2212.3157
$ws.Range("J7").Value = 2099
$ws.Range("K7").Value = 2212.3157
$ws.Range("L7").Value = 2099
$ws.Range("M7").Value = -2100.3157
$ws.Range("N7").Value = -2323
# Row 40
$ws.Range("H40").Value = 2890.4
$ws.Range("I40").Value = 2888
$ws.Range("K40").Value = 2888
$ws.Range("M40").Value = -2752
# Row 126
$ws.Range("H126").Value = 2188.7083
$ws.Range("I126").Value = 2212.3157
$ws.Range("J126").Value = 2099
$ws.Range("K126").Value = 6636.9471
$ws.Range("L126").Value = 6297
$ws.Range("M126").Value = -4166.9471
$ws.Range("N126").Value = -11237
# Row 136
$ws.Range("H136").Value = 65365.656
$ws.Range("I136").Value = 37078.605
$ws.Range("J136").Value = 263375
$ws.Range("K136").Value = 111235.815
$ws.Range("L136").Value = 790125
$ws.Range("M136").Value = -108685.815
$ws.Range("N136").Value = -795225

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 2001.3846
$ws.Range("I122").Value = 1562.5
$ws.Range("K122").Value = 4687.5
$ws.Range("M122").Value = -2237.5
# Row 126
$ws.Range("H126").Value = 1241.0526
$ws.Range("I126").Value = 1007.2727
$ws.Range("J126").Value = 1562.5
$ws.Range("K126").Value = 3021.8181
$ws.Range("L126").Value = 4687.5
$ws.Range("M126").Value = -551.8181
$ws.Range("N126").Value = -9627.5
# Row 132
$ws.Range("H132").Value = 55551.46
$ws.Range("I132").Value = 42733.168
$ws.Range("J132").Value = 79216
$ws.Range("K132").Value = 128199.504
$ws.Range("L132").Value = 237648
$ws.Range("M132").Value = -125669.504
$ws.Range("N132").Value = -242708
# Row 136
$ws.Range("H136").Value = 40122.67
$ws.Range("I136").Value = 28366.865
$ws.Range("J136").Value = 69120.336
$ws.Range("K136").Value = 85100.595
$ws.Range("L136").Value = 207361.008
$ws.Range("M136").Value = -82550.595
$ws.Range("N136").Value = -212461.008

# ---- Removal in BSM row 13 ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("N13").ClearContents()
